$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet/workbook title (reflects the new "through" date)
$ws.Name = "Through 2022-10-17"

# Update the October row label
$ws.Range("A11").Value = "October (through 10-17)"

# Row 10 (September) - 2022 column (I) count
$ws.Range("I10").Value = 146

# Row 11 (October) - per-year counts
$ws.Range("B11").Value = 14
$ws.Range("E11").Value = 45
$ws.Range("G11").Value = 82
$ws.Range("H11").Value = 106
$ws.Range("I11").Value = 56

# Row 12 (Total) - per-year counts
$ws.Range("B12").Value = 240
$ws.Range("E12").Value = 593
$ws.Range("G12").Value = 983
$ws.Range("H12").Value = 1353
$ws.Range("I12").Value = 1334
